$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, volume-change percentages).
# None of these strings are ever numeric-looking, so a normal .Value
# assignment keeps them stored as text, same as in the source file.
$textUpdates = @(
    @{ Cell = 'E2'; Value = '  -0.67%  ' },
    @{ Cell = 'E3'; Value = '  -1.34%  ' },
    @{ Cell = 'E5'; Value = '  -3.27%  ' },
    @{ Cell = 'E6'; Value = '  -2.59%  ' },
    @{ Cell = 'E7'; Value = '  -0.49%  ' },
    @{ Cell = 'E8'; Value = '  -3.71%  ' },
    @{ Cell = 'E9'; Value = '  -2.23%  ' },
    @{ Cell = 'E10'; Value = '  -2.10%  ' },
    @{ Cell = 'E11'; Value = '  -1.65%  ' },
    @{ Cell = 'E12'; Value = '  -1.39%  ' },
    @{ Cell = 'E13'; Value = '  -2.10%  ' },
    @{ Cell = 'E14'; Value = '  -3.82%  ' },
    @{ Cell = 'E15'; Value = '  -0.19%  ' },
    @{ Cell = 'E16'; Value = '  -5.18%  ' },
    @{ Cell = 'E17'; Value = '  -0.71%  ' },
    @{ Cell = 'E18'; Value = '  -0.53%  ' },
    @{ Cell = 'E19'; Value = '  -3.22%  ' },
    @{ Cell = 'E20'; Value = '  -2.07%  ' },
    @{ Cell = 'E21'; Value = '  -4.54%  ' },
    @{ Cell = 'E22'; Value = '  -1.20%  ' },
    @{ Cell = 'E23'; Value = '  -0.66%  ' },
    @{ Cell = 'E24'; Value = '  +0.45%  ' },
    @{ Cell = 'E25'; Value = '  -2.66%  ' },
    @{ Cell = 'E26'; Value = '  -4.45%  ' },
    @{ Cell = 'E27'; Value = '  +0.56%  ' },
    @{ Cell = 'E28'; Value = '  -0.63%  ' },
    @{ Cell = 'E29'; Value = '  -1.34%  ' },
    @{ Cell = 'E30'; Value = '  -3.52%  ' },
    @{ Cell = 'E32'; Value = '  -4.47%  ' },
    @{ Cell = 'E33'; Value = '  -2.68%  ' },
    @{ Cell = 'E34'; Value = '  -3.25%  ' },
    @{ Cell = 'E35'; Value = '  -1.71%  ' },
    @{ Cell = 'E36'; Value = '  -0.11%  ' },
    @{ Cell = 'E37'; Value = '  -1.34%  ' },
    @{ Cell = 'E38'; Value = '  +0.12%  ' },
    @{ Cell = 'E39'; Value = '  -0.47%  ' },
    @{ Cell = 'E40'; Value = '  -2.20%  ' },
    @{ Cell = 'E41'; Value = '  -1.19%  ' },
    @{ Cell = 'E42'; Value = '  -0.96%  ' },
    @{ Cell = 'E43'; Value = '  -1.02%  ' },
    @{ Cell = 'E44'; Value = '  -1.12%  ' },
    @{ Cell = 'E45'; Value = '  -3.71%  ' },
    @{ Cell = 'E46'; Value = '  +0.26%  ' },
    @{ Cell = 'B47'; Value = 'EnergySwap' },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Cell = 'E47'; Value = '  -1.09%  ' },
    @{ Cell = 'B48'; Value = 'Cronos' },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Cell = 'E48'; Value = '  -0.87%  ' },
    @{ Cell = 'B49'; Value = 'Mantle' },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' },
    @{ Cell = 'E49'; Value = '  -1.18%  ' },
    @{ Cell = 'B50'; Value = 'Aptos' },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Cell = 'E50'; Value = '  -1.16%  ' },
    @{ Cell = 'B51'; Value = 'Algorand' },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Cell = 'E51'; Value = '  -0.38%  ' }
)
foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Price updates. These strings look numeric (e.g. "21.21" or "26.127.36")
# so Excel would silently coerce a plain .Value assignment into a Number.
# The source workbook stores them as text, so force the Text number format
# before assigning, then drop back to the default "Normal" cell style (no
# explicit number format) so the saved cell looks just like the original.
$priceUpdates = @(
    @{ Cell = 'D2'; Value = '26.127.36' },
    @{ Cell = 'D3'; Value = '1.667.72' },
    @{ Cell = 'D5'; Value = '210.57' },
    @{ Cell = 'D6'; Value = '0.5256' },
    @{ Cell = 'D8'; Value = '0.2630' },
    @{ Cell = 'D9'; Value = '0.06301' },
    @{ Cell = 'D10'; Value = '21.21' },
    @{ Cell = 'D11'; Value = '0.07541' },
    @{ Cell = 'D12'; Value = '1.674.74' },
    @{ Cell = 'D14'; Value = '0.5564' },
    @{ Cell = 'D15'; Value = '66.78' },
    @{ Cell = 'D16'; Value = '0.000007945' },
    @{ Cell = 'D17'; Value = '26.162.26' },
    @{ Cell = 'D19'; Value = '4.748' },
    @{ Cell = 'D20'; Value = '186.60' },
    @{ Cell = 'D22'; Value = '6.180' },
    @{ Cell = 'D24'; Value = '149.73' },
    @{ Cell = 'D25'; Value = '0.1252' },
    @{ Cell = 'D27'; Value = '15.96' },
    @{ Cell = 'D29'; Value = '1.354' },
    @{ Cell = 'D31'; Value = '3.512' },
    @{ Cell = 'D33'; Value = '1.631' },
    @{ Cell = 'D35'; Value = '0.6060' },
    @{ Cell = 'D37'; Value = '2.731' },
    @{ Cell = 'D39'; Value = '1.104.84' },
    @{ Cell = 'D41'; Value = '0.8709' },
    @{ Cell = 'D43'; Value = '100.13' },
    @{ Cell = 'D44'; Value = '1.822.10' },
    @{ Cell = 'D45'; Value = '55.49' },
    @{ Cell = 'D47'; Value = '8.058' },
    @{ Cell = 'D48'; Value = '0.05237' },
    @{ Cell = 'D49'; Value = '0.4248' },
    @{ Cell = 'D50'; Value = '5.976' },
    @{ Cell = 'D51'; Value = '0.09708' }
)
foreach ($u in $priceUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}

Write-Output "Applied $($textUpdates.Count + $priceUpdates.Count) cell updates"